$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 45658

# Row 4
$ws.Range("A4").Value = 196685
$ws.Range("B4").Value = 98.53
$ws.Range("C4").Value = 19380225

# Row 5
$ws.Range("A5").Value = 2815
$ws.Range("B5").Value = 236.5
$ws.Range("C5").Value = 665754

# Row 6
$ws.Range("A6").Value = 73
$ws.Range("B6").Value = 1181.07
$ws.Range("C6").Value = 86218

# Row 7
$ws.Range("A7").Value = 199573
$ws.Range("B7").Value = 100.88
$ws.Range("C7").Value = 20132197

# Row 9
$ws.Range("A9").Value = 39862
$ws.Range("B9").Value = 155
$ws.Range("C9").Value = 6178610

# Row 10
$ws.Range("A10").Value = 80127
$ws.Range("B10").Value = 155
$ws.Range("C10").Value = 12419685

# Row 11
$ws.Range("A11").Value = 119989
$ws.Range("C11").Value = 18598295

# Row 13
$ws.Range("A13").Value = 256633
$ws.Range("B13").Value = 220
$ws.Range("C13").Value = 56459260

# Row 14
$ws.Range("A14").Value = 959
$ws.Range("B14").Value = 220
$ws.Range("C14").Value = 188080

# Row 15 - B15 cleared (value removed entirely)
$ws.Range("B15").ClearContents()

# Row 16
$ws.Range("A16").Value = 257592
$ws.Range("C16").Value = 56647340

# Row 19 - C19 cleared (value removed entirely)
$ws.Range("C19").ClearContents()

# Row 21
$ws.Range("A21").Value = 577154
$ws.Range("B21").Value = 165.26
$ws.Range("C21").Value = 95377832

# Row 22
$ws.Range("A22").Value = 377581
$ws.Range("B22").Value = 199.28
$ws.Range("C22").Value = 75245635

# Row 23
$ws.Range("A23").Value = 9510
$ws.Range("B23").Value = 239.16
$ws.Range("C23").Value = 2274381

# Update selection to F23
$ws.Range("F23").Select()
